# "Generate Report for Archive"
# The localization-status report is regenerated: every "Ready for handoff"
# status cell becomes "In Translation", and the (now narrower) Status
# columns are resized to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status text wherever it currently reads "Ready for handoff".
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value     = "In Translation"
$dede.Range("C2:C4").Value     = "In Translation"

# The Status columns now hold shorter text - shrink them to match.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth     = 12.5
$dede.Columns.Item(3).ColumnWidth     = 12.5
